$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("insurance_account")

# New test cases: acc101 (row 103) .. acc104 (row 106)
# Column layout (header row 2):
#   A Key, B LocDed6All, C LocLimit6All, D LayerAttachment, E LayerLimit,
#   F LayerParticipation, G AccDed6All, H AccMinDed6All, I AccMaxDed6All,
#   J AccLimit6All, K AccParticipation, L Supported, M STATUS

$ws.Range("A103").Value = "acc101"
$ws.Range("D103").Value = "$"
$ws.Range("E103").Value = "$"
$ws.Range("F103").Value = "%"
$ws.Range("K103").Value = "%"
$ws.Range("M103").Value = "complete"

$ws.Range("A104").Value = "acc102"
$ws.Range("K104").Value = "%"
$ws.Range("M104").Value = "complete"

$ws.Range("A105").Value = "acc103"
$ws.Range("B105").Value = "$"
$ws.Range("D105").Value = "$"
$ws.Range("E105").Value = "$"
$ws.Range("F105").Value = "%"
$ws.Range("K105").Value = "%"
$ws.Range("M105").Value = "complete"

$ws.Range("A106").Value = "acc104"
$ws.Range("B106").Value = "$"
$ws.Range("K106").Value = "%"
$ws.Range("M106").Value = "complete"
